$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Akaki"
$ws.Range("B2").Value = "Tsereteli"
$ws.Range("D2").Value = "sdgsfg"
$ws.Range("C2").Value = "fsdsdf"

$ws.Range("C2").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("D2").VerticalAlignment = -4160

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("E2").Select() | Out-Null
